$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1033.7
$ws.Range("I19").Value = 1586.25
$ws.Range("J19").Value = 665.3333
$ws.Range("K19").Value = 1586.25
$ws.Range("L19").Value = 665.3333
$ws.Range("M19").Value = -1411.25
$ws.Range("N19").Value = -1015.3333

# Row 64
$ws.Range("H64").Value = 2957.9167
$ws.Range("I64").Value = 2900
$ws.Range("J64").Value = 2977.2222
$ws.Range("K64").Value = 2900
$ws.Range("L64").Value = 2977.2222
$ws.Range("M64").Value = -2652
$ws.Range("N64").Value = -3473.2222

# Row 67
$ws.Range("H67").Value = 2957.9167
$ws.Range("I67").Value = 2900
$ws.Range("J67").Value = 2977.2222
$ws.Range("K67").Value = 2900
$ws.Range("L67").Value = 2977.2222
$ws.Range("M67").Value = -2042
$ws.Range("N67").Value = -4693.2222

# Row 70
$ws.Range("H70").Value = 1126.3572
$ws.Range("I70").Value = 1201.7273
$ws.Range("J70").Value = 850
$ws.Range("K70").Value = 3605.1819
$ws.Range("L70").Value = 2550
$ws.Range("M70").Value = -3335.1819
$ws.Range("N70").Value = -3090

# Row 73
$ws.Range("H73").Value = 1126.3572
$ws.Range("I73").Value = 1201.7273
$ws.Range("J73").Value = 850
$ws.Range("K73").Value = 3605.1819
$ws.Range("L73").Value = 2550
$ws.Range("M73").Value = -2669.1819
$ws.Range("N73").Value = -4422

# Row 74 (N74 removed)
$ws.Range("H74").Value = 3123.625
$ws.Range("I74").Value = 3123.625
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3123.625
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2187.625
$ws.Range("N74").ClearContents()

# Row 77 (N77 removed)
$ws.Range("H77").Value = 3123.625
$ws.Range("I77").Value = 3123.625
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 15618.125
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -10938.125
$ws.Range("N77").ClearContents()

# Row 103 (N103 added)
$ws.Range("H103").Value = 540.8333
$ws.Range("I103").Value = 745
$ws.Range("J103").Value = 500
$ws.Range("K103").Value = 2235
$ws.Range("L103").Value = 1500
$ws.Range("M103").Value = -1649
$ws.Range("N103").Value = -2672

# Row 129
$ws.Range("H129").Value = 2368.324
$ws.Range("J129").Value = 1044.5933
$ws.Range("L129").Value = 3133.7799
$ws.Range("N129").Value = -13133.7799

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 143676.58
$ws.Range("I45").Value = 200792.4
$ws.Range("J45").Value = 887
$ws.Range("K45").Value = 200792.4
$ws.Range("L45").Value = 887
$ws.Range("M45").Value = -200415.4
$ws.Range("N45").Value = -1641

# Row 102
$ws.Range("H102").Value = 127595
$ws.Range("I102").Value = 202116
$ws.Range("J102").Value = 3393.3333
$ws.Range("K102").Value = 202116
$ws.Range("L102").Value = 3393.3333
$ws.Range("M102").Value = -200494
$ws.Range("N102").Value = -6637.3333

# Row 122 (N122 removed)
$ws.Range("H122").Value = 1221.2222
$ws.Range("I122").Value = 1221.2222
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3663.6666
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1213.6666
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 511.42105
$ws.Range("I94").Value = 473.86667
$ws.Range("J94").Value = 652.25
$ws.Range("K94").Value = 473.86667
$ws.Range("L94").Value = 652.25
$ws.Range("M94").Value = -22.86667
$ws.Range("N94").Value = -1554.25

# Row 99
$ws.Range("H99").Value = 2281.6667
$ws.Range("I99").Value = 1895
$ws.Range("J99").Value = 2475
$ws.Range("K99").Value = 1895
$ws.Range("L99").Value = 2475
$ws.Range("M99").Value = -397
$ws.Range("N99").Value = -5471

# Row 107
$ws.Range("H107").Value = 37054468
$ws.Range("I107").Value = 58850560
$ws.Range("J107").Value = 1108
$ws.Range("K107").Value = 58850560
$ws.Range("L107").Value = 1108
$ws.Range("M107").Value = -58848640
$ws.Range("N107").Value = -4948

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 24146.64
$ws.Range("I31").Value = 1605.85
$ws.Range("J31").Value = 34392.453
$ws.Range("K31").Value = 1605.85
$ws.Range("L31").Value = 34392.453
$ws.Range("M31").Value = -1310.85
$ws.Range("N31").Value = -34982.453

# Row 34
$ws.Range("H34").Value = 24146.64
$ws.Range("I34").Value = 1605.85
$ws.Range("J34").Value = 34392.453
$ws.Range("K34").Value = 1605.85
$ws.Range("L34").Value = 34392.453
$ws.Range("M34").Value = -1403.85
$ws.Range("N34").Value = -34796.453

# Row 99
$ws.Range("H99").Value = 24194.8
$ws.Range("I99").Value = 5240
$ws.Range("J99").Value = 100014
$ws.Range("K99").Value = 5240
$ws.Range("L99").Value = 100014
$ws.Range("M99").Value = -3742
$ws.Range("N99").Value = -103010

# Row 126
$ws.Range("H126").Value = 24194.8
$ws.Range("I126").Value = 5240
$ws.Range("J126").Value = 100014
$ws.Range("K126").Value = 15720
$ws.Range("L126").Value = 300042
$ws.Range("M126").Value = -13250
$ws.Range("N126").Value = -304982

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 800.47
$ws.Range("J131").Value = 806.60205
$ws.Range("L131").Value = 2419.80615
$ws.Range("N131").Value = -12499.80615

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 4105.364
$ws.Range("I102").Value = 2646.6667
$ws.Range("J102").Value = 5855.8
$ws.Range("K102").Value = 2646.6667
$ws.Range("L102").Value = 5855.8
$ws.Range("M102").Value = -1024.6667
$ws.Range("N102").Value = -9099.799999999999

# Row 111
$ws.Range("H111").Value = 19073
$ws.Range("J111").Value = 19073
$ws.Range("L111").Value = 19073
$ws.Range("N111").Value = -25207

# Row 112
$ws.Range("H112").Value = 40799
$ws.Range("J112").Value = 40799
$ws.Range("L112").Value = 40799
$ws.Range("N112").Value = -43015

# Row 113
$ws.Range("H113").Value = 1506.2
$ws.Range("I113").Value = 940.75
$ws.Range("K113").Value = 940.75
$ws.Range("M113").Value = 1229.25

# Row 122 (N122 removed)
$ws.Range("H122").Value = 687.13043
$ws.Range("I122").Value = 687.13043
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2061.39129
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 388.60871
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 207818.14
$ws.Range("J2").Value = 7333.222
$ws.Range("L2").Value = 7333.222
$ws.Range("N2").Value = -7557.222

# Row 7
$ws.Range("H7").Value = 9074.25
$ws.Range("I7").Value = 13668
$ws.Range("J7").Value = 6318
$ws.Range("K7").Value = 13668
$ws.Range("L7").Value = 6318
$ws.Range("M7").Value = -13556
$ws.Range("N7").Value = -6542

# Row 40
$ws.Range("H40").Value = 57454.555
$ws.Range("I40").Value = 144029.72
$ws.Range("J40").Value = 2361.2727
$ws.Range("K40").Value = 144029.72
$ws.Range("L40").Value = 2361.2727
$ws.Range("M40").Value = -143893.72
$ws.Range("N40").Value = -2633.2727

# Row 61
$ws.Range("H61").Value = 1600.4584
$ws.Range("I61").Value = 1396
$ws.Range("J61").Value = 1941.2222
$ws.Range("K61").Value = 1396
$ws.Range("L61").Value = 1941.2222
$ws.Range("M61").Value = -1194
$ws.Range("N61").Value = -2345.2222

# Row 82
$ws.Range("H82").Value = 2048
$ws.Range("I82").Value = 1500
$ws.Range("J82").Value = 2596
$ws.Range("K82").Value = 1500
$ws.Range("L82").Value = 2596
$ws.Range("M82").Value = -1139
$ws.Range("N82").Value = -3318

# Row 85
$ws.Range("H85").Value = 2048
$ws.Range("I85").Value = 1500
$ws.Range("J85").Value = 2596
$ws.Range("K85").Value = 1500
$ws.Range("L85").Value = 2596
$ws.Range("M85").Value = -252
$ws.Range("N85").Value = -5092

# Row 93
$ws.Range("H93").Value = 1655.963
$ws.Range("I93").Value = 1742.421
$ws.Range("J93").Value = 1450.625
$ws.Range("K93").Value = 1742.421
$ws.Range("L93").Value = 1450.625
$ws.Range("M93").Value = -494.421
$ws.Range("N93").Value = -3946.625

# Row 113
$ws.Range("H113").Value = 1600.4584
$ws.Range("I113").Value = 1396
$ws.Range("J113").Value = 1941.2222
$ws.Range("K113").Value = 1396
$ws.Range("L113").Value = 1941.2222
$ws.Range("M113").Value = 774
$ws.Range("N113").Value = -6281.2222

# Row 122
$ws.Range("H122").Value = 2434.3684
$ws.Range("I122").Value = 2342.1538
$ws.Range("K122").Value = 7026.4614
$ws.Range("M122").Value = -4576.4614

# Row 126
$ws.Range("H126").Value = 9074.25
$ws.Range("I126").Value = 13668
$ws.Range("J126").Value = 6318
$ws.Range("K126").Value = 41004
$ws.Range("L126").Value = 18954
$ws.Range("M126").Value = -38534
$ws.Range("N126").Value = -23894

# Row 132
$ws.Range("H132").Value = 3504.6562
$ws.Range("I132").Value = 3315.5518
$ws.Range("K132").Value = 9946.6554
$ws.Range("M132").Value = -7416.6554

# Row 136
$ws.Range("H136").Value = 1713.875
$ws.Range("I136").Value = 1538.579
$ws.Range("J136").Value = 2380
$ws.Range("K136").Value = 4615.737
$ws.Range("L136").Value = 7140
$ws.Range("M136").Value = -2065.737
$ws.Range("N136").Value = -12240

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2272
$ws.Range("I122").Value = 2080.8
$ws.Range("K122").Value = 6242.400000000001
$ws.Range("M122").Value = -3792.400000000001
